$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Risk Probablity" column is inserted at C, which shifts the old
# C:O columns (Risk Level.. through the Legend table) one column to the
# right (D:P) -- this reproduces the legend table move from K:N to L:O
# and the header shuffle in one step.
$ws.Columns("C:C").Insert()

# --- Row 1 (headers) ---
$ws.Range("C1").Value = "Risk Probablity"
$ws.Range("G1").Value = "Follow up plan"
$ws.Range("L1").Value = "Legend"

# --- Row 2 ---
$ws.Range("A2").Value = "Hardware"
$ws.Range("B2").Value = "Battery Dies"
$ws.Range("C2").Value = "Med"
$ws.Range("D2").Value = "Low"
$ws.Range("G2").Value = "Charge battery when the battery is 30% or lower"
$ws.Range("L2").Value = "Types of Risks"
$ws.Range("N2").Value = "Risk Level"

# --- Row 3 ---
$ws.Range("A3").Value = "Hardware"
$ws.Range("B3").Value = "Battery Fries"
$ws.Range("C3").Value = "Low"
$ws.Range("D3").Value = "Critial"
$ws.Range("F3").Value = "2+ Days"
$ws.Range("G3").Value = "Speak to teacher"
$ws.Range("L3").Value = "Functionality"
$ws.Range("M3").Value = "Causes system to not function as required"
$ws.Range("N3").Value = "Low"
$ws.Range("O3").Value = "Will not destroy functionality"

# --- Row 4 ---
$ws.Range("A4").Value = "Hardware"
$ws.Range("B4").Value = "Battery Charger stops working"
$ws.Range("C4").Value = "Low"
$ws.Range("D4").Value = "High"
$ws.Range("F4").Value = "2+ Days"
$ws.Range("G4").Value = "Speak to teacher and barrow battery charger"
$ws.Range("L4").Value = "Bug"
$ws.Range("M4").Value = "Creates a minor issue"
$ws.Range("N4").Value = "Med"
$ws.Range("O4").Value = "Creates intermittent issues or unknowns"

# --- Row 5 ---
$ws.Range("A5").Value = "Hardware"
$ws.Range("B5").Value = "Snickerdoodle Fries"
$ws.Range("C5").Value = "Low"
$ws.Range("D5").Value = "Critial"
$ws.Range("F5").Value = "2+ Days"
$ws.Range("G5").Value = "Speak to teacher"
$ws.Range("N5").Value = "High"
$ws.Range("O5").Value = "Creates an issue to a large part of system"

# --- Row 6 ---
$ws.Range("A6").Value = "Hardware"
$ws.Range("B6").Value = "Motors burn out"
$ws.Range("C6").Value = "Low"
$ws.Range("D6").Value = "Critial"
$ws.Range("F6").Value = "2+ Days"
$ws.Range("G6").Value = "Speak to teacher"
$ws.Range("L6").Value = "Question"
$ws.Range("M6").Value = "Poses a question of not unknown"
$ws.Range("N6").Value = "Crit"
$ws.Range("O6").Value = "Breaks System"

# --- Row 7 ---
$ws.Range("A7").Value = "Hardware"
$ws.Range("B7").Value = "Motor gears are damaged"
$ws.Range("C7").Value = "Low"
$ws.Range("D7").Value = "Critial"
$ws.Range("F7").Value = "2+ Days"
$ws.Range("G7").Value = "Speak to teacher"

# --- Row 8 ---
$ws.Range("A8").Value = "Hardware"
$ws.Range("B8").Value = "Cameras are not available"
$ws.Range("C8").Value = "Med"
$ws.Range("D8").Value = "Low"
$ws.Range("E8").Value = 1000
$ws.Range("F8").Value = 2
$ws.Range("G8").Value = "Create a plan that doesn't use cameras"

# --- Row 9 ---
$ws.Range("A9").Value = "Hardware"
$ws.Range("B9").Value = "Sensor doesn't function"
$ws.Range("C9").Value = "Low"
$ws.Range("D9").Value = "Critial"
$ws.Range("F9").Value = "2+ Days"
$ws.Range("G9").Value = "Speak to teacher"

# --- Row 10 ---
$ws.Range("A10").Value = "Hardware"
$ws.Range("B10").Value = "Serial Cable not working"
$ws.Range("C10").Value = "Low"
$ws.Range("D10").Value = "Med"
$ws.Range("G10").Value = "We can swap it out"

# --- Row 11 ---
$ws.Range("A11").Value = "Environment"
$ws.Range("B11").Value = "Different floors"
$ws.Range("C11").Value = "High"
$ws.Range("D11").Value = "Med"
$ws.Range("G11").Value = "Create coeffients to handle different surfaces and/or test with other surfaces"

# --- Row 12 ---
$ws.Range("A12").Value = "Software"
$ws.Range("B12").Value = "Losing Data"
$ws.Range("C12").Value = "Med"
$ws.Range("D12").Value = "Med"
$ws.Range("G12").Value = "Backup data more often"

# --- Row 13 ---
$ws.Range("A13").Value = "System"
$ws.Range("B13").Value = "Software & Hardware Imcompability"

# --- Row 14 ---
$ws.Range("A14").Value = "Software"
$ws.Range("B14").Value = "IMU giving poor data"

# --- Row 15 ---
$ws.Range("A15").Value = "Hardware"
$ws.Range("B15").Value = "Motors connected in reverse"

# Resize the used columns to fit their new content
$ws.Columns("A:O").AutoFit()

# Restore the active selection
$ws.Range("F17").Select()
